$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.712.51'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.873.20'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '332.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.57%  '
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4737'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3956'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.69'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08025'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.022'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.81'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.876.46'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.955'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.142'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.007'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001048'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.31%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.16'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06641'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.61%  '
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.767.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.491'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.02'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.302'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.110.42'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.102'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.566'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '122.45'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9704'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09564'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.453'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.636'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.287'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06108'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02262'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.225'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.186'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.16%  '
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5984'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('E43').Value = '  +3.32%  '
$ws.Range('E44').Value = '  +0.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.252'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5689'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.30'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.408'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.932'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06816'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '112.51'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.43%  '
